$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (A:AY) get rotated in content among two independent groups of
# rows, i.e. row N ends up holding the full set of values that used to live in
# row Sigma(N) (everything else in the sheet stays untouched).
#
#   Group 1: 2 <- 3, 3 <- 23, 23 <- 2
#   Group 2: 10 <- 24, 11 <- 10, 12 <- 11, 13 <- 12, 14 <- 13, 15 <- 14,
#            16 <- 15, 17 <- 16, 18 <- 17, 19 <- 18, 20 <- 19, 21 <- 20,
#            22 <- 21, 24 <- 22

$lastCol = "AY"

# Snapshot the full row contents (columns A:AY) for every row that is a
# source in the rotation, before any writes happen, so that later writes
# never clobber data that is still needed.
$srcRows = @(2, 3, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24)
$snapshot = @{}
foreach ($r in $srcRows) {
    $snapshot[$r] = $ws.Range("A$r`:$lastCol$r").Value2
}

# Destination row -> source row (content that should be copied into it)
$sigma = @{
    2  = 3
    3  = 23
    23 = 2
    10 = 24
    11 = 10
    12 = 11
    13 = 12
    14 = 13
    15 = 14
    16 = 15
    17 = 16
    18 = 17
    19 = 18
    20 = 19
    21 = 20
    22 = 21
    24 = 22
}

foreach ($destRow in $sigma.Keys) {
    # Columns Y and AA store free-text dates (e.g. "2019-05-08") as plain
    # strings in the source workbook. Force the destination cells to Text
    # format first so that Excel does not silently reinterpret the string
    # as a real date serial number when the value is written back in.
    $ws.Range("Y$destRow").NumberFormat = "@"
    $ws.Range("AA$destRow").NumberFormat = "@"

    $ws.Range("A$destRow`:$lastCol$destRow").Value2 = $snapshot[$sigma[$destRow]]
}

Write-Host "Row rotation complete"
